# Update "想去人数" (column F) values across sheets per latest scrape.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 781
$ws1.Range("F3").Value = 970
$ws1.Range("F4").Value = 752
$ws1.Range("F5").Value = 857
$ws1.Range("F6").Value = 420
$ws1.Range("F7").Value = 648
$ws1.Range("F8").Value = 144
$ws1.Range("F9").Value = 1245
$ws1.Range("F10").Value = 669
$ws1.Range("F11").Value = 399
$ws1.Range("F12").Value = 527
$ws1.Range("F13").Value = 173
$ws1.Range("F15").Value = 714
$ws1.Range("F17").Value = 380
$ws1.Range("F18").Value = 362
$ws1.Range("F20").Value = 567
$ws1.Range("F21").Value = 116
$ws1.Range("F22").Value = 606
$ws1.Range("F23").Value = 31
$ws1.Range("F24").Value = 867

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 25
$ws2.Range("F8").Value = 232
$ws2.Range("F11").Value = 25

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 781
$ws4.Range("F6").Value = 970
$ws4.Range("F7").Value = 752
$ws4.Range("F8").Value = 857
$ws4.Range("F9").Value = 420
$ws4.Range("F10").Value = 648
$ws4.Range("F11").Value = 144
$ws4.Range("F12").Value = 1245
$ws4.Range("F13").Value = 669
$ws4.Range("F15").Value = 25
$ws4.Range("F16").Value = 399
$ws4.Range("F17").Value = 527
$ws4.Range("F19").Value = 173
$ws4.Range("F21").Value = 714
$ws4.Range("F24").Value = 380
$ws4.Range("F25").Value = 362
$ws4.Range("F27").Value = 232
$ws4.Range("F29").Value = 567
$ws4.Range("F31").Value = 25
$ws4.Range("F34").Value = 116
$ws4.Range("F35").Value = 606
$ws4.Range("F36").Value = 31
$ws4.Range("F37").Value = 867

$wb.Save()
